$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 6 ("logistic") to make room for "scaledTanh",
# shifting logistic..sin down by one row (rows 6-11 -> 7-12).
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with the "scaledTanh" activation function results.
$ws.Range("A6").Value = "scaledTanh"
$ws.Range("B6").Value = 45.471400000000003
$ws.Range("C6").Value = 41.372500000000002
$ws.Range("D6").Value = 17.931999999999999
$ws.Range("E6").Value = 13.785600000000001
$ws.Range("F6").Value = 11.2369
$ws.Range("G6").Value = $null
$ws.Range("H6").Value = $null
$ws.Range("I6").Value = $null

# Update selection / view to match the author's final state.
$ws.Range("D9").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
